$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5890.76
$ws.Range("J17").Value = 1929.3469
$ws.Range("L17").Value = 5788.0407
$ws.Range("N17").Value = -6124.0407
$ws.Range("H33").Value = 666.6
$ws.Range("I33").Value = 547.8461
$ws.Range("K33").Value = 547.8461
$ws.Range("M33").Value = -318.8461
$ws.Range("H58").Value = 22694.389
$ws.Range("I58").Value = 401
$ws.Range("J58").Value = 25227.727
$ws.Range("K58").Value = 1203
$ws.Range("L58").Value = 75683.181
$ws.Range("M58").Value = -1053
$ws.Range("N58").Value = -75983.181
$ws.Range("H70").Value = 3696.3635
$ws.Range("I70").Value = 1200
$ws.Range("J70").Value = 3946
$ws.Range("K70").Value = 3600
$ws.Range("L70").Value = 11838
$ws.Range("M70").Value = -3330
$ws.Range("N70").Value = -12378
$ws.Range("H73").Value = 3696.3635
$ws.Range("I73").Value = 1200
$ws.Range("J73").Value = 3946
$ws.Range("K73").Value = 3600
$ws.Range("L73").Value = 11838
$ws.Range("M73").Value = -2664
$ws.Range("N73").Value = -13710
$ws.Range("H74").Value = 4146.6665
$ws.Range("J74").Value = 4146.6665
$ws.Range("L74").Value = 4146.6665
$ws.Range("N74").Value = -6018.6665
$ws.Range("H77").Value = 4146.6665
$ws.Range("J77").Value = 4146.6665
$ws.Range("L77").Value = 20733.3325
$ws.Range("N77").Value = -30093.3325
$ws.Range("H82").Value = 4980.6665
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 4980.6665
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H96").Value = 790.5294
$ws.Range("I96").Value = 738.1
$ws.Range("K96").Value = 2214.3
$ws.Range("M96").Value = -841.3000000000002
$ws.Range("H135").Value = 1024.8889
$ws.Range("I135").Value = 754.42426
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 6789.81834
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -4254.81834
$ws.Range("N135").Value = -41070
$ws.Range("H138").Value = 4135.519
$ws.Range("I138").Value = 1091.8572
$ws.Range("J138").Value = 4791.077
$ws.Range("K138").Value = 3275.5716
$ws.Range("L138").Value = 14373.231
$ws.Range("M138").Value = 1864.4284
$ws.Range("N138").Value = -24653.231
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 73335
$ws.Range("J10").Value = 73335
$ws.Range("L10").Value = 73335
$ws.Range("N10").Value = -73675
$ws.Range("H61").Value = 2926.2368
$ws.Range("I61").Value = 1372.6364
$ws.Range("K61").Value = 1372.6364
$ws.Range("M61").Value = -1160.6364
$ws.Range("H74").Value = 878.5405
$ws.Range("I74").Value = 806.1739
$ws.Range("J74").Value = 997.4286
$ws.Range("K74").Value = 806.1739
$ws.Range("L74").Value = 997.4286
$ws.Range("M74").Value = 67.8261
$ws.Range("N74").Value = -2745.4286
$ws.Range("H77").Value = 878.5405
$ws.Range("I77").Value = 806.1739
$ws.Range("J77").Value = 997.4286
$ws.Range("K77").Value = 4030.8695
$ws.Range("L77").Value = 4987.143
$ws.Range("M77").Value = 337.1305000000002
$ws.Range("N77").Value = -13723.143
$ws.Range("H131").Value = 23200
$ws.Range("J131").Value = 23200
$ws.Range("L131").Value = 23200
$ws.Range("N131").Value = -33280
$ws.Range("H132").Value = 14086618
$ws.Range("I132").Value = 21277748
$ws.Range("K132").Value = 63833244
$ws.Range("M132").Value = -63830714
$ws.Range("H136").Value = 2926.2368
$ws.Range("I136").Value = 1372.6364
$ws.Range("K136").Value = 4117.9092
$ws.Range("M136").Value = -1567.9092
$ws.Range("H141").Value = 31200
$ws.Range("J141").Value = 31200
$ws.Range("L141").Value = 31200
$ws.Range("N141").Value = -41560
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4549.6
$ws.Range("I134").Value = 3299.4
$ws.Range("K134").Value = 9898.200000000001
$ws.Range("M134").Value = -7363.200000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5003
$ws.Range("I31").Value = 3389.8
$ws.Range("J31").Value = 8229.4
$ws.Range("K31").Value = 3389.8
$ws.Range("L31").Value = 8229.4
$ws.Range("M31").Value = -3094.8
$ws.Range("N31").Value = -8819.4
$ws.Range("H34").Value = 5003
$ws.Range("I34").Value = 3389.8
$ws.Range("J34").Value = 8229.4
$ws.Range("K34").Value = 3389.8
$ws.Range("L34").Value = 8229.4
$ws.Range("M34").Value = -3187.8
$ws.Range("N34").Value = -8633.4
$ws.Range("H48").Value = 51593.285
$ws.Range("J48").Value = 51593.285
$ws.Range("L48").Value = 51593.285
$ws.Range("N48").Value = -52545.285
$ws.Range("H50").Value = 13676.8
$ws.Range("J50").Value = 13676.8
$ws.Range("L50").Value = 13676.8
$ws.Range("N50").Value = -14926.8
$ws.Range("H51").Value = 14079.6
$ws.Range("J51").Value = 14079.6
$ws.Range("L51").Value = 14079.6
$ws.Range("N51").Value = -15551.6
$ws.Range("H59").Value = 10140
$ws.Range("J59").Value = 10140
$ws.Range("L59").Value = 10140
$ws.Range("N59").Value = -12430
$ws.Range("H60").Value = 14467.667
$ws.Range("J60").Value = 14701.5
$ws.Range("L60").Value = 14701.5
$ws.Range("N60").Value = -15723.5
$ws.Range("H61").Value = 14079.6
$ws.Range("J61").Value = 14079.6
$ws.Range("L61").Value = 14079.6
$ws.Range("N61").Value = -14775.6
$ws.Range("H68").Value = 30899.285
$ws.Range("J68").Value = 30899.285
$ws.Range("L68").Value = 30899.285
$ws.Range("N68").Value = -32397.285
$ws.Range("H71").Value = 30899.285
$ws.Range("J71").Value = 30899.285
$ws.Range("L71").Value = 92697.855
$ws.Range("N71").Value = -100185.855
$ws.Range("H134").Value = 2360.238
$ws.Range("I134").Value = 1059.1818
$ws.Range("J134").Value = 3791.4
$ws.Range("K134").Value = 3177.5454
$ws.Range("L134").Value = 11374.2
$ws.Range("M134").Value = -642.5454
$ws.Range("N134").Value = -16444.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 11699.25
$ws.Range("I3").Value = 783.3333
$ws.Range("J3").Value = 18248.8
$ws.Range("K3").Value = 783.3333
$ws.Range("L3").Value = 18248.8
$ws.Range("M3").Value = -667.3333
$ws.Range("N3").Value = -18480.8
$ws.Range("H7").Value = 5718428.5
$ws.Range("J7").Value = 5000000
$ws.Range("L7").Value = 5000000
$ws.Range("N7").Value = -5000224
$ws.Range("H8").Value = 5718428.5
$ws.Range("J8").Value = 5000000
$ws.Range("L8").Value = 5000000
$ws.Range("N8").Value = -5000278
$ws.Range("H10").Value = 19200.8
$ws.Range("I10").Value = 8000
$ws.Range("J10").Value = 36002
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 36002
$ws.Range("M10").Value = -7831
$ws.Range("N10").Value = -36340
$ws.Range("H80").Value = 4247.5
$ws.Range("I80").Value = 4557.5
$ws.Range("J80").Value = 3937.5
$ws.Range("K80").Value = 4557.5
$ws.Range("L80").Value = 3937.5
$ws.Range("M80").Value = -3559.5
$ws.Range("N80").Value = -5933.5
$ws.Range("H83").Value = 4247.5
$ws.Range("I83").Value = 4557.5
$ws.Range("J83").Value = 3937.5
$ws.Range("K83").Value = 22787.5
$ws.Range("L83").Value = 19687.5
$ws.Range("M83").Value = -17795.5
$ws.Range("N83").Value = -29671.5
$ws.Range("H126").Value = 3467.8823
$ws.Range("I126").Value = 2603
$ws.Range("J126").Value = 3734
$ws.Range("K126").Value = 7809
$ws.Range("L126").Value = 11202
$ws.Range("M126").Value = -5339
$ws.Range("N126").Value = -16142
$ws.Range("H132").Value = 2992.8086
$ws.Range("I132").Value = 2589.1614
$ws.Range("J132").Value = 3774.875
$ws.Range("K132").Value = 7767.4842
$ws.Range("L132").Value = 11324.625
$ws.Range("M132").Value = -5237.4842
$ws.Range("N132").Value = -16384.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1953.2069
$ws.Range("I68").Value = 1051.6666
$ws.Range("J68").Value = 6280.6
$ws.Range("K68").Value = 1051.6666
$ws.Range("L68").Value = 6280.6
$ws.Range("M68").Value = -302.6666
$ws.Range("N68").Value = -7778.6
$ws.Range("H71").Value = 1953.2069
$ws.Range("I71").Value = 1051.6666
$ws.Range("J71").Value = 6280.6
$ws.Range("K71").Value = 5258.333000000001
$ws.Range("L71").Value = 31403
$ws.Range("M71").Value = -1514.333000000001
$ws.Range("N71").Value = -38891
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H131").Value = 19933.334
$ws.Range("J131").Value = 19933.334
$ws.Range("L131").Value = 19933.334
$ws.Range("N131").Value = -30013.334
$ws.Range("H132").Value = 3298.1667
$ws.Range("I132").Value = 1984.4667
$ws.Range("J132").Value = 4611.8667
$ws.Range("K132").Value = 5953.4001
$ws.Range("L132").Value = 13835.6001
$ws.Range("M132").Value = -3423.4001
$ws.Range("N132").Value = -18895.6001
$ws.Range("H135").Value = 29775.572
$ws.Range("J135").Value = 29775.572
$ws.Range("L135").Value = 29775.572
$ws.Range("N135").Value = -39915.572
$ws.Range("H137").Value = 29833.334
$ws.Range("J137").Value = 29833.334
$ws.Range("L137").Value = 29833.334
$ws.Range("N137").Value = -40033.334
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H110").Value = 36500
$ws.Range("J110").Value = 36500
$ws.Range("L110").Value = 36500
$ws.Range("N110").Value = -44680

Write-Output "done"